$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.52 = 25799.73 pesos`n✅ 25799.73 pesos = 6.5 = 975.32 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the tasas worksheet figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 153.49
$wsTasas.Range("O10").Value = 3960
$wsTasas.Range("N12").Value = 3967.9
$wsTasas.Range("O12").Value = 150
